# Kayıt silindi: 11292947
# Remove the deleted record's row from both the master "Kayitlar" list
# and its filtered per-department view ("Merkez İlçe"); EntireRow.Delete
# shifts all following rows up by one, matching the diff.

$wb = $excel.ActiveWorkbook

$wsKayitlar = $wb.Worksheets.Item("Kayitlar")
$wsKayitlar.Rows.Item(1131).Delete()

$wsMerkez = $wb.Worksheets.Item("Merkez İlçe")
$wsMerkez.Rows.Item(592).Delete()
